$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (including the date/time number format) from the row above,
# then set the new date value for the new row's date cell.
$ws.Range("A31").Copy($ws.Range("A32"))
$ws.Range("A32").Value = 45936

# New quote values for the added row (stored as text, matching the sheet's
# existing convention of comma-decimal text values).
$ws.Range("B32").Value = "21,5922"
$ws.Range("C32").Value = "15,3073"
$ws.Range("D32").Value = "15,4736"
$ws.Range("E32").Value = "15,4736"
